# Apply the edit described by the diff:
# A new data row is inserted right before the current row 53 ("Hortaliza, Vega Central
# Mapocho de Santiago - Arveja Verde" sheet), pushing all existing rows 53..113 down by
# one (53->54, 54->55, ..., 112->113, 113->114). The new row 53 receives a brand-new
# price observation.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row above the current row 53. Excel shifts rows 53:113 down to 54:114,
# preserving their existing content/formatting (matches the diff's row-shift pattern).
$ws.Rows.Item(53).Insert()

# Populate the newly inserted row 53 with the new record's values.
$ws.Range("A53").Value = 9
$ws.Range("B53").Value = "Vega Central Mapocho de Santiago"
$ws.Range("C53").Value = "Metropolitana"
$ws.Range("D53").Value = 44671
$ws.Range("E53").Value = 13
$ws.Range("F53").Value = 100112022
$ws.Range("G53").Value = "Arveja Verde"
$ws.Range("H53").Value = "Sin especificar"
$ws.Range("I53").Value = "Primera"
$ws.Range("J53").Value = 25
$ws.Range("K53").Value = 26000
$ws.Range("L53").Value = 27000
$ws.Range("M53").Value = 26520
$ws.Range("N53").Value = '$/saco 25 kilos'
$ws.Range("O53").Value = "Carahue"
$ws.Range("P53").Value = 1061
$ws.Range("Q53").Value = 25
$ws.Range("R53").Value = "Hortaliza"
